{"js": "// Update the benchmark results table. The table has a single column and\n// each row holds one metric. A handful of rows got new/condensed values:\n// some summary numbers changed, and three rows that used to hold a whole\n// tab-separated breakdown (process id + per-phase timings + \"100.0\") were\n// collapsed down to just their leading count.\nconst tables = context.document.body.tables;\ntables.load('items');\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst updates = [\n  [0, '0M'],\n  [1, '0M'],\n  [2, '0M'],\n  [3, '32'],\n  [4, '0.00002'],\n  [8, '0.00004'],\n  [11, '0.00115'],\n  [43, '100'],\n  [44, '0'],\n  [45, '102'],\n];\n\nfor (const [rowIndex, newValue] of updates) {\n  const cell = table.getCell(rowIndex, 0);\n  cell.value = newValue;\n}\n\nawait context.sync();\n", "ps1": "# Update the benchmark results table. The table has a single column and\n# each row holds one metric. A handful of rows got new/condensed values:\n# some summary numbers changed, and three rows that used to hold a whole\n# tab-separated breakdown (process id + per-phase timings + \"100.0\") were\n# collapsed down to just their leading count.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$updates = @{\n    1  = \"0M\"\n    2  = \"0M\"\n    3  = \"0M\"\n    4  = \"32\"\n    5  = \"0.00002\"\n    9  = \"0.00004\"\n    12 = \"0.00115\"\n    44 = \"100\"\n    45 = \"0\"\n    46 = \"102\"\n}\n\nforeach ($rowIndex in $updates.Keys) {\n    $cell = $t.Cell($rowIndex, 1)\n    $cell.Range.Text = $updates[$rowIndex]\n}\n"}
